# feat: add the import function and fix some bugs
#
# 1) Sheet "2024-12": insert two new columns (类型/账户) after 摘要, fill them
#    in for every existing row, append two more transaction rows, and update
#    the 支出/结余 (now shifted to G/H) values to their new post-import totals.
# 2) Add a new sheet "2023-01" (after "2024-12") holding the imported ledger
#    rows for January 2023, with the same header layout.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "2024-12"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("2024-12")

# Make room for the two new columns (类型, 账户) right after 摘要 (C). This
# shifts the old D:F (收入/支出/结余) to F:H.
$ws1.Range("D1:E1").EntireColumn.Insert()

$ws1.Range("D1").Value = "类型"
$ws1.Range("E1").Value = "账户"

# Existing rows (2-6) all describe "dinner" / "food" purchases; fill in the
# new 类型/账户 columns and refresh the running balance in H now that the
# ledger has been reconciled against the imported accounts.
$ws1.Range("D2").Value = "food"
$ws1.Range("E2").Value = "alipay"
$ws1.Range("H2").Value = 1029.6

$ws1.Range("D3").Value = "food"
$ws1.Range("E3").Value = "alipay"
$ws1.Range("H3").Value = 1100

$ws1.Range("D4").Value = "food"
$ws1.Range("E4").Value = "alipay"
$ws1.Range("H4").Value = 1009.5999999999999

$ws1.Range("D5").Value = "food"
$ws1.Range("E5").Value = "alipay"
$ws1.Range("H5").Value = 1019.5999999999999

$ws1.Range("D6").Value = "food"
$ws1.Range("E6").Value = "wechat"
$ws1.Range("G6").Value = 30
$ws1.Range("H6").Value = 1070

# Two new rows imported from the wechat/alipay statements.
$ws1.Range("A7").Value = 12
$ws1.Range("B7").Value = 8
$ws1.Range("C7").Value = "dinner"
$ws1.Range("D7").Value = "food"
$ws1.Range("E7").Value = "wechat"
$ws1.Range("G7").Value = 20.4
$ws1.Range("H7").Value = 1049.6

$ws1.Range("A8").Value = 12
$ws1.Range("B8").Value = 8
$ws1.Range("C8").Value = "dinner"
$ws1.Range("D8").Value = "food"
$ws1.Range("E8").Value = "alipay"
$ws1.Range("G8").Value = 10
$ws1.Range("H8").Value = 1039.6

# ---------------------------------------------------------------------------
# New sheet "2023-01" (imported historical ledger)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws2.Name = "2023-01"

$ws2.Range("A1").Value = "月"
$ws2.Range("B1").Value = "日"
$ws2.Range("C1").Value = "摘要"
$ws2.Range("D1").Value = "类型"
$ws2.Range("E1").Value = "账户"
$ws2.Range("F1").Value = "收入"
$ws2.Range("G1").Value = "支出"
$ws2.Range("H1").Value = "结余"

$ws2.Range("A2").Value = 1
$ws2.Range("B2").Value = 8
$ws2.Range("C2").Value = "dinner"
$ws2.Range("D2").Value = "food"
$ws2.Range("E2").Value = "alipay"
$ws2.Range("G2").Value = 23
$ws2.Range("H2").Value = 986.5999999999999

$ws2.Range("A3").Value = 1
$ws2.Range("B3").Value = 8
$ws2.Range("C3").Value = "dinner"
$ws2.Range("D3").Value = "food"
$ws2.Range("E3").Value = "alipay"
$ws2.Range("G3").Value = 23
$ws2.Range("H3").Value = 963.5999999999999

$ws2.Range("A4").Value = 1
$ws2.Range("B4").Value = 8
$ws2.Range("C4").Value = "dinner"
$ws2.Range("D4").Value = "food"
$ws2.Range("E4").Value = "alipay"
$ws2.Range("G4").Value = 23
$ws2.Range("H4").Value = 940.5999999999999

$ws2.Range("A5").Value = 1
$ws2.Range("B5").Value = 8
$ws2.Range("C5").Value = "dinner"
$ws2.Range("D5").Value = "food"
$ws2.Range("E5").Value = "alipay"
$ws2.Range("G5").Value = 23
$ws2.Range("H5").Value = 917.5999999999999

# Keep "2024-12" as the active/selected tab, as it was before the edit.
$ws1.Activate()
